$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows (in descending order so earlier row numbers stay valid)
# corresponds to removed accounts: 005870700, 002973105, 004216504, 005666419,
# 004222784, 004231509, 004382374
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()
